$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows 2-11 (columns B:K) with latest run values
$ws.Cells.Item(2, 2).Value = 0.32577036678551236
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0

$ws.Cells.Item(3, 2).Value = 0.30843618033632086
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = -0.0009265441771828323
$ws.Cells.Item(3, 5).Value = [double]"4.213297079274692e-05"
$ws.Cells.Item(3, 6).Value = -0.0005735477300171508
$ws.Cells.Item(3, 7).Value = 0.0002720611532660484
$ws.Cells.Item(3, 8).Value = [double]"-8.72095251145644e-06"
$ws.Cells.Item(3, 9).Value = -0.0003691483579295538
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0.004414442729769208

$ws.Cells.Item(4, 2).Value = 0.3057621424924377
$ws.Cells.Item(4, 3).Value = -0.0005146833720973999
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = [double]"1.4045034288427827e-05"
$ws.Cells.Item(4, 6).Value = [double]"-3.0295275694206335e-06"
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = [double]"-1.7575599740610715e-06"
$ws.Cells.Item(4, 9).Value = -0.0007349632850170942
$ws.Cells.Item(4, 10).Value = [double]"-5.3446802259467764e-05"
$ws.Cells.Item(4, 11).Value = -0.00046511758164330397

$ws.Cells.Item(5, 2).Value = 0.2918930825326562
$ws.Cells.Item(5, 3).Value = 0.0029878728614583278
$ws.Cells.Item(5, 4).Value = -0.001891167284993457
$ws.Cells.Item(5, 5).Value = [double]"5.581158208826706e-06"
$ws.Cells.Item(5, 6).Value = -0.0002462232858239164
$ws.Cells.Item(5, 7).Value = -0.001065090961855057
$ws.Cells.Item(5, 8).Value = [double]"-6.704567252877974e-05"
$ws.Cells.Item(5, 9).Value = [double]"-6.917556286548641e-05"
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = -0.0005009637092014407

$ws.Cells.Item(6, 2).Value = 0.36645096516291703
$ws.Cells.Item(6, 3).Value = 0.015506121168330361
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = [double]"-1.9073547315090085e-05"
$ws.Cells.Item(6, 6).Value = [double]"-4.6296280922476205e-05"
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = [double]"-6.219135003561905e-05"
$ws.Cells.Item(6, 9).Value = -0.0019803415867228856
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = [double]"4.04403304109513e-05"

$ws.Cells.Item(7, 2).Value = 0.36045042016366957
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = -0.0018529242856495343
$ws.Cells.Item(7, 5).Value = -0.00019084837197465746
$ws.Cells.Item(7, 6).Value = -0.001180379824683229
$ws.Cells.Item(7, 7).Value = 0.0005704111613147494
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0.0001223419975827628
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0.0017358371720212484

$ws.Cells.Item(8, 2).Value = 0.24921101423408276
$ws.Cells.Item(8, 3).Value = -0.015974001822534164
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = [double]"3.8301988688003535e-05"
$ws.Cells.Item(8, 6).Value = -0.00029548013311802364
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = [double]"2.391500383187952e-05"
$ws.Cells.Item(8, 9).Value = 0.00025942179423668426
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0.0002795735389726772

$ws.Cells.Item(9, 2).Value = 0.23998846500863788
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0.00016971053206171156
$ws.Cells.Item(9, 5).Value = -0.0017611994298077438
$ws.Cells.Item(9, 6).Value = -0.005829534744179096
$ws.Cells.Item(9, 7).Value = 0.0004479659521581601
$ws.Cells.Item(9, 8).Value = -0.00018511762898431968
$ws.Cells.Item(9, 9).Value = [double]"-4.470951208916138e-05"
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0.0007002779205462495

$ws.Cells.Item(10, 2).Value = 0.4232943953622883
$ws.Cells.Item(10, 3).Value = 0.038175018279621796
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = [double]"-3.0201578051124305e-05"
$ws.Cells.Item(10, 6).Value = -0.0002941827035763721
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = [double]"1.1569427451533095e-05"
$ws.Cells.Item(10, 9).Value = -0.00025573727888316154
$ws.Cells.Item(10, 10).Value = -0.002031163010122183
$ws.Cells.Item(10, 11).Value = 0.00020391618029075298

$ws.Cells.Item(11, 2).Value = 0.4170504736592897
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0.0006524386395050954
$ws.Cells.Item(11, 5).Value = 0.00027643857934656286
$ws.Cells.Item(11, 6).Value = -0.0012250954763661094
$ws.Cells.Item(11, 7).Value = 0.0009818170403482662
$ws.Cells.Item(11, 8).Value = [double]"2.6784021978306805e-05"
$ws.Cells.Item(11, 9).Value = -0.0022090327842295247
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0.008667609888168615

# Add new row 12 for 2025-08-30 (force text so it is not auto-converted to a date serial)
$ws.Cells.Item(12, 1).Value = "'2025-08-30"
$ws.Cells.Item(12, 1).Style = "Normal"
$ws.Cells.Item(12, 2).Value = 0.34857488907790923
$ws.Cells.Item(12, 3).Value = -0.02798273831566987
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = [double]"4.1289423430865034e-05"
$ws.Cells.Item(12, 6).Value = [double]"6.74841952477652e-06"
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = [double]"-5.854738601282607e-06"
$ws.Cells.Item(12, 9).Value = -0.00105733310251161
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = [double]"6.238915264544787e-05"

# Adjust column widths for E and H to match latest layout
$ws.Columns.Item(5).ColumnWidth = 15.77734375
$ws.Columns.Item(8).ColumnWidth = 15.77734375
